# Remove unused parameter for retirement home
# Row 10 on the "common_parameters" sheet holds
# perc_single_over90yo_in_retirement_home, which is no longer used.
# Delete the entire row and let everything below shift up by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = $ws.Rows.Item(10)
$row.Select()
$row.EntireRow.Delete()

# Land the selection where Excel would leave it after a row delete:
# the row that slid up into the deleted row's place.
$ws.Range("A10").Select()
